$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New values for columns B (Belegdatum), C (Lieferant), D (Profitcenter)
# for rows 2-7. Column B merges the old date-day with the old month text,
# column C takes the old Profitcenter letter, and column D gets a new
# project code.
$data = @(
    @{ Row = 2; B = "05.Jan"; C = "A"; D = "P1" },
    @{ Row = 3; B = "06.Jan"; C = "B"; D = "P3" },
    @{ Row = 4; B = "07.Jan"; C = "A"; D = "P1" },
    @{ Row = 5; B = "09.Jan"; C = "B"; D = "P2" },
    @{ Row = 6; B = "10.Jan"; C = "A"; D = "P3" },
    @{ Row = 7; B = "11.Jan"; C = "A"; D = "P2" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
}
